$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, border, centered) from H1 into I1:J1, then set header text
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill I2:J41 data values (columns I=9, J=10)
$data = @(
    @(2, 7, 11),
    @(3, 9, 9),
    @(4, 7, 8),
    @(5, 8, 8),
    @(6, 8, 9),
    @(7, 6, 8),
    @(8, 5, 7),
    @(9, 6, 7),
    @(10, 9, 9),
    @(11, 2, 5),
    @(12, 6, 6),
    @(13, 9, 9),
    @(14, 6, 6),
    @(15, 5, 7),
    @(16, 2, 5),
    @(17, 3, 8),
    @(18, 1, 5),
    @(19, 1, 6),
    @(20, 1, 5),
    @(21, 2, 7),
    @(22, 1, 5),
    @(23, 1, 6),
    @(24, 1, 5),
    @(25, 1, 5),
    @(26, 1, 6),
    @(27, 1, 6),
    @(28, 1, 5),
    @(29, 4, 5),
    @(30, 1, 4),
    @(31, 5, 7),
    @(32, 8, 8),
    @(33, 6, 8),
    @(34, 7, 8),
    @(35, 8, 8),
    @(36, 1, 3),
    @(37, 6, 9),
    @(38, 1, 3),
    @(39, 1, 3),
    @(40, 6, 7),
    @(41, 1, 1)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}

Write-Output "Done"